# Apply forecast pivot update:
# - Row 2: Project Owner -> "-"; monthly values Aug'25..Oct'26 (H:V) 10000 -> 5333.33
# - Row 3: Project Owner -> "-"
# - Row 4: Project Owner -> "-"; E4 20000 -> 16000; I4 40000 -> 32000; M4 40000 -> 32000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"
$ws.Range("H2:V2").Value = 5333.33

# Row 3
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 16000
$ws.Range("I4").Value = 32000
$ws.Range("M4").Value = 32000
